$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 865.3333
$ws.Range("I4").Value = 798.25
$ws.Range("K4").Value = 798.25
$ws.Range("M4").Value = -684.25
$ws.Range("H9").Value = 180.33333
$ws.Range("I9").Value = 187.375
$ws.Range("K9").Value = 187.375
$ws.Range("M9").Value = -18.375
$ws.Range("H18").Value = 959.5
$ws.Range("I18").Value = 959.5
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 959.5
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -675.5
$ws.Range("N18").Value = ""
$ws.Range("H28").Value = 3584
$ws.Range("J28").Value = 6248.125
$ws.Range("L28").Value = 6248.125
$ws.Range("N28").Value = -7218.125
$ws.Range("H39").Value = 170.16667
$ws.Range("I39").Value = 181
$ws.Range("K39").Value = 543
$ws.Range("M39").Value = -247
$ws.Range("H55").Value = 358.22223
$ws.Range("I55").Value = 368.8889
$ws.Range("J55").Value = 347.55554
$ws.Range("K55").Value = 368.8889
$ws.Range("L55").Value = 347.55554
$ws.Range("M55").Value = -154.8889
$ws.Range("N55").Value = -775.5555400000001
$ws.Range("H62").Value = 8650.1
$ws.Range("I62").Value = 2498.3333
$ws.Range("K62").Value = 2498.3333
$ws.Range("M62").Value = -1874.3333
$ws.Range("H65").Value = 8650.1
$ws.Range("I65").Value = 2498.3333
$ws.Range("K65").Value = 12491.6665
$ws.Range("M65").Value = -9371.6665
$ws.Range("H86").Value = 9579.2
$ws.Range("I86").Value = 9299
$ws.Range("K86").Value = 9299
$ws.Range("M86").Value = -8176
$ws.Range("H89").Value = 9579.2
$ws.Range("I89").Value = 9299
$ws.Range("K89").Value = 46495
$ws.Range("M89").Value = -40879
$ws.Range("H116").Value = 3982
$ws.Range("J116").Value = 3960
$ws.Range("L116").Value = 3960
$ws.Range("N116").Value = -10844
$ws.Range("H125").Value = 3191.25
$ws.Range("J125").Value = 3000
$ws.Range("L125").Value = 27000
$ws.Range("N125").Value = -31920
$ws.Range("H127").Value = 878.5
$ws.Range("I127").Value = 878.5
$ws.Range("K127").Value = 2635.5
$ws.Range("M127").Value = 2324.5
$ws.Range("H129").Value = 1146.6666
$ws.Range("I129").Value = 915
$ws.Range("K129").Value = 2745
$ws.Range("M129").Value = 2255
$ws.Range("H132").Value = 9268.111
$ws.Range("I132").Value = 10773.333
$ws.Range("K132").Value = 32319.999
$ws.Range("M132").Value = -29789.999
$ws.Range("H140").Value = 44000
$ws.Range("J140").Value = 44000
$ws.Range("L140").Value = 44000
$ws.Range("N140").Value = -54360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 77.57143
$ws.Range("I5").Value = 68.6
$ws.Range("K5").Value = 68.6
$ws.Range("M5").Value = 43.40000000000001
$ws.Range("H39").Value = 6424.75
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").Value = ""
$ws.Range("H95").Value = 9000
$ws.Range("J95").Value = 9000
$ws.Range("L95").Value = 9000
$ws.Range("N95").Value = -14492
$ws.Range("H96").Value = 20000000
$ws.Range("J96").Value = 20000000
$ws.Range("L96").Value = 20000000
$ws.Range("N96").Value = -20005492
$ws.Range("H102").Value = 3117.2942
$ws.Range("I102").Value = 1445.6923
$ws.Range("J102").Value = 8550
$ws.Range("K102").Value = 1445.6923
$ws.Range("L102").Value = 8550
$ws.Range("M102").Value = 176.3077000000001
$ws.Range("N102").Value = -11794
$ws.Range("H122").Value = 1233.1666
$ws.Range("I122").Value = 1233.1666
$ws.Range("K122").Value = 3699.4998
$ws.Range("M122").Value = -1249.4998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 77.57143
$ws.Range("I4").Value = 68.6
$ws.Range("K4").Value = 68.6
$ws.Range("M4").Value = 46.40000000000001
$ws.Range("H11").Value = 875.75
$ws.Range("I11").Value = 875.75
$ws.Range("K11").Value = 875.75
$ws.Range("M11").Value = -735.75
$ws.Range("H86").Value = 3049.0715
$ws.Range("J86").Value = 5458.0835
$ws.Range("L86").Value = 5458.0835
$ws.Range("N86").Value = -7704.0835
$ws.Range("H89").Value = 3049.0715
$ws.Range("J89").Value = 5458.0835
$ws.Range("L89").Value = 27290.4175
$ws.Range("N89").Value = -38522.4175
$ws.Range("H94").Value = 296.33334
$ws.Range("I94").Value = 294.5
$ws.Range("J94").Value = 300
$ws.Range("K94").Value = 294.5
$ws.Range("L94").Value = 300
$ws.Range("M94").Value = 156.5
$ws.Range("N94").Value = -1202
$ws.Range("H105").Value = 1455.4
$ws.Range("I105").Value = 1455.4
$ws.Range("K105").Value = 1455.4
$ws.Range("M105").Value = 291.5999999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 4510.3477
$ws.Range("I7").Value = 5370.3687
$ws.Range("J7").Value = 425.25
$ws.Range("K7").Value = 5370.3687
$ws.Range("L7").Value = 425.25
$ws.Range("M7").Value = -5257.3687
$ws.Range("N7").Value = -651.25
$ws.Range("H22").Value = 2149.5
$ws.Range("I22").Value = 1849.5
$ws.Range("J22").Value = 2449.5
$ws.Range("K22").Value = 1849.5
$ws.Range("L22").Value = 2449.5
$ws.Range("M22").Value = -1499.5
$ws.Range("N22").Value = -3149.5
$ws.Range("H122").Value = 743.36365
$ws.Range("I122").Value = 467.7
$ws.Range("K122").Value = 1403.1
$ws.Range("M122").Value = 1046.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 15464.857
$ws.Range("I50").Value = 17209
$ws.Range("J50").Value = 5000
$ws.Range("K50").Value = 51627
$ws.Range("L50").Value = 15000
$ws.Range("M50").Value = -51146
$ws.Range("N50").Value = -15962
$ws.Range("H53").Value = 15464.857
$ws.Range("I53").Value = 17209
$ws.Range("J53").Value = 5000
$ws.Range("K53").Value = 51627
$ws.Range("L53").Value = 15000
$ws.Range("M53").Value = -51146
$ws.Range("N53").Value = -15962
$ws.Range("H103").Value = 343.1
$ws.Range("I103").Value = 352.33334
$ws.Range("J103").Value = 329.25
$ws.Range("K103").Value = 1057.00002
$ws.Range("L103").Value = 987.75
$ws.Range("M103").Value = -178.0000199999999
$ws.Range("N103").Value = -2745.75
$ws.Range("H132").Value = 4874.75
$ws.Range("I132").Value = 4499
$ws.Range("J132").Value = 4928.4287
$ws.Range("K132").Value = 40491
$ws.Range("L132").Value = 44355.85830000001
$ws.Range("M132").Value = -37961
$ws.Range("N132").Value = -49415.85830000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 345.77274
$ws.Range("I2").Value = 194.6923
$ws.Range("J2").Value = 564
$ws.Range("K2").Value = 194.6923
$ws.Range("L2").Value = 564
$ws.Range("M2").Value = -81.69229999999999
$ws.Range("N2").Value = -790
$ws.Range("H70").Value = 9000
$ws.Range("I70").Value = 9000
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 9000
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -8730
$ws.Range("N70").Value = ""
$ws.Range("H73").Value = 9000
$ws.Range("I73").Value = 9000
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 9000
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -8064
$ws.Range("N73").Value = ""
$ws.Range("H80").Value = 2800.2
$ws.Range("J80").Value = 2502.75
$ws.Range("L80").Value = 2502.75
$ws.Range("N80").Value = -4498.75
$ws.Range("H83").Value = 2800.2
$ws.Range("J83").Value = 2502.75
$ws.Range("L83").Value = 12513.75
$ws.Range("N83").Value = -22497.75
$ws.Range("H102").Value = 3371.7144
$ws.Range("I102").Value = 3267
$ws.Range("K102").Value = 3267
$ws.Range("M102").Value = -1645
$ws.Range("H113").Value = 4858.6
$ws.Range("I113").Value = 1431.3334
$ws.Range("J113").Value = 9999.5
$ws.Range("K113").Value = 1431.3334
$ws.Range("L113").Value = 9999.5
$ws.Range("M113").Value = 738.6666
$ws.Range("N113").Value = -14339.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 6104.1577
$ws.Range("I46").Value = 4497.5
$ws.Range("J46").Value = 6532.6
$ws.Range("K46").Value = 4497.5
$ws.Range("L46").Value = 6532.6
$ws.Range("M46").Value = -4309.5
$ws.Range("N46").Value = -6908.6
$ws.Range("H68").Value = 7964.2856
$ws.Range("I68").Value = 6250
$ws.Range("J68").Value = 8250
$ws.Range("K68").Value = 6250
$ws.Range("L68").Value = 8250
$ws.Range("M68").Value = -5501
$ws.Range("N68").Value = -9748
$ws.Range("H71").Value = 7964.2856
$ws.Range("I71").Value = 6250
$ws.Range("J71").Value = 8250
$ws.Range("K71").Value = 31250
$ws.Range("L71").Value = 41250
$ws.Range("M71").Value = -27506
$ws.Range("N71").Value = -48738

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 13350833
$ws.Range("I5").Value = 15025000
$ws.Range("J5").Value = 10002500
$ws.Range("K5").Value = 15025000
$ws.Range("L5").Value = 10002500
$ws.Range("M5").Value = -15024888
$ws.Range("N5").Value = -10002724
$ws.Range("H62").Value = 8846.462
$ws.Range("J62").Value = 9625.25
$ws.Range("L62").Value = 9625.25
$ws.Range("N62").Value = -10873.25
$ws.Range("H65").Value = 8846.462
$ws.Range("J65").Value = 9625.25
$ws.Range("L65").Value = 48126.25
$ws.Range("N65").Value = -54366.25
$ws.Range("H122").Value = 2680.8667
$ws.Range("I122").Value = 1646.2727
$ws.Range("K122").Value = 4938.8181
$ws.Range("M122").Value = -2488.8181

